$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("✅ 1000 Bs = 5.36 = 21371.05 pesos", "✅ 1000 Bs = 5.13 = 20456.88 pesos")
$text = $text.Replace("✅ 21371.05 pesos = 5.33 = 953.73 Bs", "✅ 20456.88 pesos = 5.13 = 942.86 Bs")
$cell.Value2 = $text

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 194.8
$wsTasas.Range("O10").Value2 = 3985
$wsTasas.Range("N12").Value2 = 3990
$wsTasas.Range("O12").Value2 = 183.9
